$wb = $excel.ActiveWorkbook

# --- Sheet1 ("added"): insert a Date column after Name, before Email ---
$ws1 = $wb.Worksheets.Item("added")

# Insert a new column before column B (Email), shifting Email/Phone/Address/Status right
$ws1.Columns.Item(2).Insert()

# Header for new Date column (style is already inherited from the Insert above)
$ws1.Cells.Item(1, 2).Value = "Date"

$dateValue = 45905.039861111109

$c2 = $ws1.Cells.Item(2, 2)
$c2.Value = $dateValue

$c3 = $ws1.Cells.Item(3, 2)
$c3.Value = $dateValue

# --- Sheet3 ("common"): rework phone column values ---
$ws3 = $wb.Worksheets.Item("common")

# "0350" is entered as quote-prefixed text (keeps leading zero, stored as General/text)
$p2 = $ws3.Cells.Item(2, 3)
$p2.Value = "'0350"

$p3 = $ws3.Cells.Item(3, 3)
$p3.Value = 3500

$p4 = $ws3.Cells.Item(4, 3)
$p4.Value = -5600

$ws3.Activate()
$ws3.Range("D4").Select() | Out-Null

# --- back to Sheet1: apply the long-date number format to the new Date values ---
$longDateFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$c2.NumberFormat = $longDateFormat
$c3.NumberFormat = $longDateFormat

$ws1.Range("A1:F3").EntireColumn.AutoFit() | Out-Null

$ws1.Activate()
$ws1.Range("B2").Select() | Out-Null
